$d = $word.ActiveDocument

# Helper: replace the text of an (italic, placeholder-red) paragraph with new
# text, dropping the red color but keeping italic - mirrors typing fresh text
# over a cleared placeholder run.
function Set-ItalicNoColor($para, [string]$text) {
    $r = $para.Range
    $r.End = $r.End - 1
    $r.Delete()
    $r.InsertAfter($text)
    $r.Font.Italic = $true
    $r.Font.ItalicBi = $true
}

# ---------------------------------------------------------------------------
# 1. "TextView" -> "MultilineText" (list of Layout view kinds)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("TextView", $true, $false, $false, $false, $false, $true, 1, $false, "MultilineText", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Fill in the "example JSON reply" placeholder paragraph with the real
#    Open Trivia DB JSON payload (drop the red placeholder color, keep italic)
# ---------------------------------------------------------------------------
$jsonReply = '{"category":"General Knowledge","type":"multiple","difficulty":"easy","question":"What geometric shape is generally used for stop signs?","correct_answer":"Octagon","incorrect_answers":["Hexagon","Circle","Triangle"]},{"category":"General Knowledge","type":"multiple","difficulty":"easy","question":"What is the shape of the toy invented by Hungarian professor Ernu0151 Rubik?","correct_answer":"Cube","incorrect_answers":["Sphere","Cylinder","Pyramid"]},{"category":"General Knowledge","type":"multiple","difficulty":"easy","question":"What machine element is located in the center of fidget spinners?","correct_answer":"Bearings","incorrect_answers":["Axles","Gears","Belts"]},{"category":"General Knowledge","type":"multiple","difficulty":"easy","question":"How many furlongs are there in a mile?","correct_answer":"Eight","incorrect_answers":["Two","Four","Six"]},{"category":"General Knowledge","type":"multiple","difficulty":"easy","question":"Earth is located in which galaxy?","correct_answer":"The Milky Way Galaxy","incorrect_answers":["The Mars Galaxy","The Galaxy Note","The Black Hole"]}]}""'

$rng = $d.Content
$rng.Find.Execute("Insert JSON reply here") | Out-Null
$jsonPara1 = $rng.Paragraphs(1)
Set-ItalicNoColor $jsonPara1 $jsonReply
